$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 23.51647972924624
$ws.Range("C2").Value = 16.07267116366772
$ws.Range("D2").Value = 6.096773559099828
$ws.Range("E2").Value = 12.32158857521578
$ws.Range("F2").Value = 48.68540974036375
$ws.Range("J2").Value = 10.44424916278694
$ws.Range("N2").Value = 19.69065503691025

$ws.Range("B3").Value = 23.0145074299656
$ws.Range("C3").Value = 15.60482678393161
$ws.Range("D3").Value = 6.101291582118657
$ws.Range("E3").Value = 12.3072696705451
$ws.Range("F3").Value = 48.38303298373083
$ws.Range("J3").Value = 10.45396233582858
$ws.Range("N3").Value = 19.75688350268317

$ws.Range("B4").Value = 22.70867302754717
$ws.Range("C4").Value = 15.31605601974642
$ws.Range("D4").Value = 6.104604796827453
$ws.Range("E4").Value = 12.30116094845537
$ws.Range("F4").Value = 48.21136809616085
$ws.Range("J4").Value = 10.46192875252664
$ws.Range("N4").Value = 19.79957821313474

$ws.Range("B5").Value = 22.58483208960049
$ws.Range("C5").Value = 15.1982017218248
$ws.Range("D5").Value = 6.106091166166886
$ws.Range("E5").Value = 12.29934718394972
$ws.Range("F5").Value = 48.14497429615949
$ws.Range("J5").Value = 10.46567755428539
$ws.Range("N5").Value = 19.81748732873605

$ws.Range("B6").Value = 22.56432167178715
$ws.Range("C6").Value = 15.17862747146446
$ws.Range("D6").Value = 6.106346221854107
$ws.Range("E6").Value = 12.29908683065355
$ws.Range("F6").Value = 48.13416582303385
$ws.Range("J6").Value = 10.46633034969076
$ws.Range("N6").Value = 19.82049197287153

$ws.Range("B7").Value = 22.70699941662847
$ws.Range("C7").Value = 15.31446704082626
$ws.Range("D7").Value = 6.104624290180243
$ws.Range("E7").Value = 12.30113375110285
$ws.Range("F7").Value = 48.21045821734291
$ws.Range("J7").Value = 10.46197727745137
$ws.Range("N7").Value = 19.79981767375295

$ws.Range("B8").Value = 23.34301977724232
$ws.Range("C8").Value = 15.91179058953543
$ws.Range("D8").Value = 6.098219861105689
$ws.Range("E8").Value = 12.31609519766948
$ws.Range("F8").Value = 48.57827243357139
$ws.Range("J8").Value = 10.44718201107494
$ws.Range("N8").Value = 19.71306935240536

$ws.Range("B9").Value = 24.60090203205927
$ws.Range("C9").Value = 17.06279902722078
$ws.Range("D9").Value = 6.089910612619495
$ws.Range("E9").Value = 12.36666690305398
$ws.Range("F9").Value = 49.40847033938044
$ws.Range("J9").Value = 10.43410329698465
$ws.Range("N9").Value = 19.55905582560059

$ws.Range("B10").Value = 25.52091696912129
$ws.Range("C10").Value = 17.88574570553721
$ws.Range("D10").Value = 6.086359240176684
$ws.Range("E10").Value = 12.41666295741432
$ws.Range("F10").Value = 50.08155185634322
$ws.Range("J10").Value = 10.43426526414422
$ws.Range("N10").Value = 19.45570354246039

$ws.Range("B11").Value = 25.93644470233327
$ws.Range("C11").Value = 18.2532840735699
$ws.Range("D11").Value = 6.085290033697714
$ws.Range("E11").Value = 12.44216427130272
$ws.Range("F11").Value = 50.40065492126667
$ws.Range("J11").Value = 10.43646987167375
$ws.Range("N11").Value = 19.41081215462908

$ws.Range("B12").Value = 26.09320034133984
$ws.Range("C12").Value = 18.39133915044106
$ws.Range("D12").Value = 6.08496302451061
$ws.Range("E12").Value = 12.45221404385866
$ws.Range("F12").Value = 50.52327279289811
$ws.Range("J12").Value = 10.43761158056322
$ws.Range("N12").Value = 19.39411844397333

$ws.Range("B13").Value = 26.05946923166853
$ws.Range("C13").Value = 18.36165862133893
$ws.Range("D13").Value = 6.085029999223004
$ws.Range("E13").Value = 12.45003222986595
$ws.Range("F13").Value = 50.49678691627503
$ws.Range("J13").Value = 10.4373520394863
$ws.Range("N13").Value = 19.39770013340602

$ws.Range("B14").Value = 25.94935393400155
$ws.Range("C14").Value = 18.26466523656779
$ws.Range("D14").Value = 6.085261573653336
$ws.Range("E14").Value = 12.44298321850112
$ws.Range("F14").Value = 50.41070753745556
$ws.Range("J14").Value = 10.43655764858526
$ws.Range("N14").Value = 19.4094326250778

$ws.Range("B15").Value = 25.88182281670311
$ws.Range("C15").Value = 18.20510365016587
$ws.Range("D15").Value = 6.085413540981613
$ws.Range("E15").Value = 12.43871655829137
$ws.Range("F15").Value = 50.35821095794891
$ws.Range("J15").Value = 10.4361110355576
$ws.Range("N15").Value = 19.41665893213049

$ws.Range("B16").Value = 25.4936866805133
$ws.Range("C16").Value = 17.86157640403691
$ws.Range("D16").Value = 6.086440050457778
$ws.Range("E16").Value = 12.41505158465917
$ws.Range("F16").Value = 50.06095113938623
$ws.Range("J16").Value = 10.43416410940215
$ws.Range("N16").Value = 19.45868006177298

$ws.Range("B17").Value = 25.25469061114692
$ws.Range("C17").Value = 17.64898006997489
$ws.Range("D17").Value = 6.087209223568594
$ws.Range("E17").Value = 12.40123782990418
$ws.Range("F17").Value = 49.88184644016667
$ws.Range("J17").Value = 10.43351588324456
$ws.Range("N17").Value = 19.485002848527

$ws.Range("B18").Value = 25.1169540925607
$ws.Range("C18").Value = 17.52606559685909
$ws.Range("D18").Value = 6.087703096274421
$ws.Range("E18").Value = 12.39355235122115
$ws.Range("F18").Value = 49.78005063558193
$ws.Range("J18").Value = 10.43334359705199
$ws.Range("N18").Value = 19.50034293039594

$ws.Range("B19").Value = 25.0702773108722
$ws.Range("C19").Value = 17.48434442296579
$ws.Range("D19").Value = 6.087879172969981
$ws.Range("E19").Value = 12.39099490791813
$ws.Range("F19").Value = 49.74579621305494
$ws.Range("J19").Value = 10.43331969342211
$ws.Range("N19").Value = 19.50557114686009

$ws.Range("B20").Value = 25.28016156330072
$ws.Range("C20").Value = 17.67167807488969
$ws.Range("D20").Value = 6.087122022324404
$ws.Range("E20").Value = 12.40268146111002
$ws.Range("F20").Value = 49.90078666864937
$ws.Range("J20").Value = 10.43356412746944
$ws.Range("N20").Value = 19.48218005131834

$ws.Range("B21").Value = 25.98171490233227
$ws.Range("C21").Value = 18.29318611028437
$ws.Range("D21").Value = 6.08519144659313
$ws.Range("E21").Value = 12.4450430496328
$ws.Range("F21").Value = 50.43594344694665
$ws.Range("J21").Value = 10.43678264925827
$ws.Range("N21").Value = 19.40597820949625

$ws.Range("B22").Value = 26.43667881039369
$ws.Range("C22").Value = 18.69277101539927
$ws.Range("D22").Value = 6.084383302320561
$ws.Range("E22").Value = 12.4750173866305
$ws.Range("F22").Value = 50.79603903029432
$ws.Range("J22").Value = 10.44067491858678
$ws.Range("N22").Value = 19.3579578101232

$ws.Range("B23").Value = 26.1942321412623
$ws.Range("C23").Value = 18.48015309726509
$ws.Range("D23").Value = 6.084773348114202
$ws.Range("E23").Value = 12.45881141894158
$ws.Range("F23").Value = 50.60292948250871
$ws.Range("J23").Value = 10.43843376107116
$ws.Range("N23").Value = 19.38342406181032

$ws.Range("B24").Value = 25.26864717504165
$ws.Range("C24").Value = 17.66141844513803
$ws.Range("D24").Value = 6.087161285030521
$ws.Range("E24").Value = 12.40202799706599
$ws.Range("F24").Value = 49.89222012923322
$ws.Range("J24").Value = 10.43354169212808
$ws.Range("N24").Value = 19.48345559393972

$ws.Range("B25").Value = 24.26059631667611
$ws.Range("C25").Value = 16.75471240777014
$ws.Range("D25").Value = 6.091707345534112
$ws.Range("E25").Value = 12.35072157389215
$ws.Range("F25").Value = 49.17253397481906
$ws.Range("J25").Value = 10.43592908034517
$ws.Range("N25").Value = 19.59899853923191

